# Refresh crypto symbol price/volume data (GitHub Actions scheduled update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'301.55"
$ws.Cells.Item(2, 5).Value = "'-4.32%"
$ws.Cells.Item(3, 4).Value = "'35.14"
$ws.Cells.Item(3, 5).Value = "'-0.33%"
$ws.Cells.Item(4, 4).Value = "'5.036"
$ws.Cells.Item(4, 5).Value = "'-1.98%"
$ws.Cells.Item(5, 4).Value = "'0.07984"
$ws.Cells.Item(5, 5).Value = "'-1.77%"
$ws.Cells.Item(6, 4).Value = "'1.894"
$ws.Cells.Item(6, 5).Value = "'-11.23%"
$ws.Cells.Item(7, 4).Value = "'7.806"
$ws.Cells.Item(7, 5).Value = "'-1.93%"
$ws.Cells.Item(8, 5).Value = "'-2.34%"
$ws.Cells.Item(10, 4).Value = "'0.9226"
$ws.Cells.Item(10, 5).Value = "'-0.76%"
$ws.Cells.Item(11, 4).Value = "'0.1266"
$ws.Cells.Item(11, 5).Value = "'25.70%"
$ws.Cells.Item(12, 5).Value = "'-0.98%"
$ws.Cells.Item(13, 4).Value = "'0.1003"
$ws.Cells.Item(13, 5).Value = "'9.61%"
$ws.Cells.Item(14, 4).Value = "'0.03528"
$ws.Cells.Item(14, 5).Value = "'-2.04%"
$ws.Cells.Item(15, 4).Value = "'0.09845"
$ws.Cells.Item(15, 5).Value = "'-0.63%"
$ws.Cells.Item(16, 4).Value = "'0.001391"
$ws.Cells.Item(16, 5).Value = "'-3.38%"
$ws.Cells.Item(17, 4).Value = "'0.005889"
$ws.Cells.Item(17, 5).Value = "'3.05%"
$ws.Cells.Item(18, 4).Value = "'3.505"
$ws.Cells.Item(18, 5).Value = "'1.11%"
$ws.Cells.Item(19, 5).Value = "'-0.32%"
$ws.Cells.Item(20, 5).Value = "'-2.88%"
$ws.Cells.Item(21, 4).Value = "'5.032"
$ws.Cells.Item(21, 5).Value = "'-1.19%"
$ws.Cells.Item(22, 5).Value = "'8.08%"
$ws.Cells.Item(23, 4).Value = "'0.04500"
$ws.Cells.Item(23, 5).Value = "'-1.57%"
$ws.Cells.Item(24, 4).Value = "'0.001214"
$ws.Cells.Item(24, 5).Value = "'-2.57%"
$ws.Cells.Item(25, 4).Value = "'0.004788"
$ws.Cells.Item(25, 5).Value = "'1.72%"
$ws.Cells.Item(26, 4).Value = "'0.0001250"
$ws.Cells.Item(26, 5).Value = "'-0.07%"
$ws.Cells.Item(27, 5).Value = "'-33.37%"
$ws.Cells.Item(39, 4).Value = "'0.01884"
$ws.Cells.Item(39, 5).Value = "'-3.66%"
$ws.Cells.Item(40, 4).Value = "'0.04720"
$ws.Cells.Item(40, 5).Value = "'-2.57%"
$ws.Cells.Item(41, 4).Value = "'0.007510"
$ws.Cells.Item(41, 5).Value = "'-2.68%"
$ws.Cells.Item(42, 4).Value = "'0.01024"
$ws.Cells.Item(42, 5).Value = "'30.48%"
$ws.Cells.Item(43, 4).Value = "'0.1324"
$ws.Cells.Item(43, 5).Value = "'-4.91%"
$ws.Cells.Item(44, 4).Value = "'0.002110"
$ws.Cells.Item(44, 5).Value = "'-2.83%"
$ws.Cells.Item(45, 4).Value = "'0.01056"
$ws.Cells.Item(45, 5).Value = "'-10.49%"
$ws.Cells.Item(46, 4).Value = "'0.00006221"
$ws.Cells.Item(46, 5).Value = "'-6.04%"
$ws.Cells.Item(47, 4).Value = "'0.00000000750"
$ws.Cells.Item(47, 5).Value = "'-0.06%"
$ws.Cells.Item(48, 5).Value = "'71.75%"
$ws.Cells.Item(49, 5).Value = "'-12.45%"
$ws.Cells.Item(50, 4).Value = "'0.00002100"
$ws.Cells.Item(50, 5).Value = "'-0.06%"
$ws.Cells.Item(51, 4).Value = "'0.0002000"
$ws.Cells.Item(51, 5).Value = "'-0.06%"
